$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 497.75
$ws.Range("I11").Value = 497.75
$ws.Range("K11").Value = 497.75
$ws.Range("M11").Value = -357.75
$ws.Range("H18").Value = 250
$ws.Range("I18").Value = 250
$ws.Range("K18").Value = 250
$ws.Range("M18").Value = 34
$ws.Range("H43").Value = 4813.2
$ws.Range("I43").Value = 2174.75
$ws.Range("J43").Value = 5772.636
$ws.Range("K43").Value = 2174.75
$ws.Range("L43").Value = 5772.636
$ws.Range("M43").Value = -2105.75
$ws.Range("N43").Value = -5910.636
$ws.Range("H112").Value = 3703.3684
$ws.Range("I112").Value = 18388.5
$ws.Range("J112").Value = 1975.7059
$ws.Range("K112").Value = 55165.5
$ws.Range("L112").Value = 5927.1177
$ws.Range("M112").Value = -54057.5
$ws.Range("N112").Value = -8143.1177
$ws.Range("H137").Value = 2505.7144
$ws.Range("I137").Value = 1700.8572
$ws.Range("K137").Value = 5102.571599999999
$ws.Range("M137").Value = -2552.571599999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 487.7
$ws.Range("I4").Value = 468.625
$ws.Range("K4").Value = 468.625
$ws.Range("M4").Value = -352.625
$ws.Range("H12").Value = 1525.75
$ws.Range("I12").Value = 1366.3334
$ws.Range("J12").Value = 2004
$ws.Range("K12").Value = 1366.3334
$ws.Range("L12").Value = 2004
$ws.Range("M12").Value = -1193.3334
$ws.Range("N12").Value = -2350
$ws.Range("H32").Value = 1985465.6
$ws.Range("I32").Value = 622.7347
$ws.Range("K32").Value = 622.7347
$ws.Range("M32").Value = -335.7347
$ws.Range("H61").Value = 3852.7334
$ws.Range("I61").Value = 2258.2
$ws.Range("J61").Value = 4650
$ws.Range("K61").Value = 2258.2
$ws.Range("L61").Value = 4650
$ws.Range("M61").Value = -2046.2
$ws.Range("N61").Value = -5074
$ws.Range("H74").Value = 1778.6666
$ws.Range("I74").Value = 1815.5714
$ws.Range("K74").Value = 1815.5714
$ws.Range("M74").Value = -941.5714
$ws.Range("H77").Value = 1778.6666
$ws.Range("I77").Value = 1815.5714
$ws.Range("K77").Value = 9077.857
$ws.Range("M77").Value = -4709.857
$ws.Range("H97").Value = 875.5238000000001
$ws.Range("I97").Value = 718.7646999999999
$ws.Range("K97").Value = 718.7646999999999
$ws.Range("M97").Value = -222.7646999999999
$ws.Range("H122").Value = 3124.75
$ws.Range("I122").Value = 2067.4285
$ws.Range("K122").Value = 6202.2855
$ws.Range("M122").Value = -3752.2855
$ws.Range("H132").Value = 6249406
$ws.Range("I132").Value = 10103010
$ws.Range("K132").Value = 30309030
$ws.Range("M132").Value = -30306500
$ws.Range("H133").Value = 172249.25
$ws.Range("J133").Value = 172249.25
$ws.Range("L133").Value = 172249.25
$ws.Range("N133").Value = -177309.25
$ws.Range("H136").Value = 3852.7334
$ws.Range("I136").Value = 2258.2
$ws.Range("J136").Value = 4650
$ws.Range("K136").Value = 6774.599999999999
$ws.Range("L136").Value = 13950
$ws.Range("M136").Value = -4224.599999999999
$ws.Range("N136").Value = -19050

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2664.6667
$ws.Range("J94").Value = 5514.8335
$ws.Range("L94").Value = 5514.8335
$ws.Range("N94").Value = -6416.8335
$ws.Range("H107").Value = 10001248
$ws.Range("I107").Value = 11112163
$ws.Range("K107").Value = 11112163
$ws.Range("M107").Value = -11110243
$ws.Range("H134").Value = 1912131.4
$ws.Range("I134").Value = 2384515.8
$ws.Range("K134").Value = 7153547.399999999
$ws.Range("M134").Value = -7151012.399999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 468.6
$ws.Range("I7").Value = 398
$ws.Range("J7").Value = 574.5
$ws.Range("K7").Value = 398
$ws.Range("L7").Value = 574.5
$ws.Range("M7").Value = -285
$ws.Range("N7").Value = -800.5
$ws.Range("H16").Value = 55560704
$ws.Range("J16").Value = 7237.5
$ws.Range("L16").Value = 7237.5
$ws.Range("N16").Value = -7811.5
$ws.Range("H113").Value = 55560704
$ws.Range("J113").Value = 7237.5
$ws.Range("L113").Value = 7237.5
$ws.Range("N113").Value = -11577.5
$ws.Range("H134").Value = 33339650
$ws.Range("I134").Value = 250001630
$ws.Range("J134").Value = 7038.3076
$ws.Range("K134").Value = 750004890
$ws.Range("L134").Value = 21114.9228
$ws.Range("M134").Value = -750002355
$ws.Range("N134").Value = -26184.9228

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 1074.7142
$ws.Range("I6").Value = 1074.7142
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3224.1426
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -3111.1426
$ws.Range("N6").ClearContents()
$ws.Range("H37").Value = 218246.28
$ws.Range("J37").Value = 218246.28
$ws.Range("L37").Value = 654738.84
$ws.Range("N37").Value = -654962.84
$ws.Range("H38").Value = 103.4
$ws.Range("J38").Value = 120.5
$ws.Range("L38").Value = 361.5
$ws.Range("N38").Value = -1055.5
$ws.Range("H129").Value = 41667004
$ws.Range("I129").Value = 449.66666
$ws.Range("J129").Value = 166666670
$ws.Range("K129").Value = 1348.99998
$ws.Range("L129").Value = 500000010
$ws.Range("M129").Value = 3651.00002
$ws.Range("N129").Value = -500010010

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1099.3214
$ws.Range("I97").Value = 1028.9259
$ws.Range("K97").Value = 1028.9259
$ws.Range("M97").Value = -532.9259
$ws.Range("H122").Value = 10223.55
$ws.Range("I122").Value = 11009.75
$ws.Range("J122").Value = 9699.416999999999
$ws.Range("K122").Value = 33029.25
$ws.Range("L122").Value = 29098.251
$ws.Range("M122").Value = -30579.25
$ws.Range("N122").Value = -33998.251
$ws.Range("H126").Value = 20007804
$ws.Range("J126").Value = 19166.5
$ws.Range("L126").Value = 57499.5
$ws.Range("N126").Value = -62439.5
$ws.Range("H132").Value = 58828092
$ws.Range("I132").Value = 100003760
$ws.Range("J132").Value = 5714
$ws.Range("K132").Value = 300011280
$ws.Range("L132").Value = 17142
$ws.Range("M132").Value = -300008750
$ws.Range("N132").Value = -22202

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4724.115
$ws.Range("I7").Value = 2928.3845
$ws.Range("J7").Value = 6519.846
$ws.Range("K7").Value = 2928.3845
$ws.Range("L7").Value = 6519.846
$ws.Range("M7").Value = -2816.3845
$ws.Range("N7").Value = -6743.846
$ws.Range("H22").Value = 1070.7142
$ws.Range("I22").Value = 1149.5
$ws.Range("J22").Value = 965.6667
$ws.Range("K22").Value = 1149.5
$ws.Range("L22").Value = 965.6667
$ws.Range("M22").Value = -854.5
$ws.Range("N22").Value = -1555.6667
$ws.Range("H27").Value = 1070.7142
$ws.Range("I27").Value = 1149.5
$ws.Range("J27").Value = 965.6667
$ws.Range("K27").Value = 1149.5
$ws.Range("L27").Value = 965.6667
$ws.Range("M27").Value = -1042.5
$ws.Range("N27").Value = -1179.6667
$ws.Range("H61").Value = 5136.4116
$ws.Range("I61").Value = 3555
$ws.Range("K61").Value = 3555
$ws.Range("M61").Value = -3353
$ws.Range("H68").Value = 2066.3333
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 2066.3333
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value = 2899.4
$ws.Range("I93").Value = 3298.2
$ws.Range("K93").Value = 3298.2
$ws.Range("M93").Value = -2050.2
$ws.Range("H100").Value = 3243.842
$ws.Range("I100").Value = 4854.375
$ws.Range("K100").Value = 4854.375
$ws.Range("M100").Value = -4313.375
$ws.Range("H113").Value = 5136.4116
$ws.Range("I113").Value = 3555
$ws.Range("K113").Value = 3555
$ws.Range("M113").Value = -1385
$ws.Range("H122").Value = 5165.5557
$ws.Range("I122").Value = 4257
$ws.Range("J122").Value = 5743.727
$ws.Range("K122").Value = 12771
$ws.Range("L122").Value = 17231.181
$ws.Range("M122").Value = -10321
$ws.Range("N122").Value = -22131.181
$ws.Range("H126").Value = 4724.115
$ws.Range("I126").Value = 2928.3845
$ws.Range("J126").Value = 6519.846
$ws.Range("K126").Value = 8785.1535
$ws.Range("L126").Value = 19559.538
$ws.Range("M126").Value = -6315.1535
$ws.Range("N126").Value = -24499.538
$ws.Range("H132").Value = 2557.6365
$ws.Range("I132").Value = 1851.1875
$ws.Range("K132").Value = 5553.5625
$ws.Range("M132").Value = -3023.5625

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17071
$ws.Range("J62").Value = 16349.4
$ws.Range("L62").Value = 16349.4
$ws.Range("N62").Value = -17597.4
$ws.Range("H65").Value = 17071
$ws.Range("J65").Value = 16349.4
$ws.Range("L65").Value = 81747
$ws.Range("N65").Value = -87987
$ws.Range("H112").Value = 27749.5
$ws.Range("J112").Value = 27749.5
$ws.Range("L112").Value = 27749.5
$ws.Range("N112").Value = -30703.5
$ws.Range("H122").Value = 2743.8125
$ws.Range("I122").Value = 2239.5
$ws.Range("K122").Value = 6718.5
$ws.Range("M122").Value = -4268.5
$ws.Range("H132").Value = 4205.607
$ws.Range("I132").Value = 2569.3684
$ws.Range("J132").Value = 7659.8887
$ws.Range("K132").Value = 7708.1052
$ws.Range("L132").Value = 22979.6661
$ws.Range("M132").Value = -5178.1052
$ws.Range("N132").Value = -28039.6661
$ws.Range("H136").Value = 20879306
$ws.Range("J136").Value = 8819.625
$ws.Range("L136").Value = 26458.875
$ws.Range("N136").Value = -31558.875
